$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 59

# Column A holds a date-formatted string ("2025/10/04") that must stay literal
# text (matching the other rows), not get auto-converted into a date serial
# number. Force text entry via NumberFormat, then strip the formatting back
# off so the cell ends up with the default (unstyled) look, same as its
# siblings.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/04"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "土"
$ws.Cells.Item($row, 3).Value = 8
$ws.Cells.Item($row, 4).Value = 201
